$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New daily data rows (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
# covering 2021-09-21 through 2021-12-08 (commit: "aggiornamento fino a 8/12").
$data = @(
  @(44460, 0, 7, 40.85205719288007),
  @(44461, 0, 7, 40.85205719288007),
  @(44462, 2, 7, 40.85205719288007),
  @(44463, 2, 7, 40.85205719288007),
  @(44464, 0, 7, 40.85205719288007),
  @(44465, 0, 5, 29.18004085205719),
  @(44466, 1, 5, 29.18004085205719),
  @(44467, 1, 6, 35.01604902246864),
  @(44468, 0, 6, 35.01604902246864),
  @(44469, 1, 5, 29.18004085205719),
  @(44470, 2, 5, 29.18004085205719),
  @(44471, 2, 7, 40.85205719288007),
  @(44472, 0, 7, 40.85205719288007),
  @(44473, 2, 8, 46.68806536329151),
  @(44474, 4, 11, 64.19608987452582),
  @(44475, 0, 11, 64.19608987452582),
  @(44476, 4, 14, 81.70411438576014),
  @(44477, 3, 15, 87.54012255617158),
  @(44478, 0, 13, 75.8681062153487),
  @(44479, 0, 13, 75.8681062153487),
  @(44480, 2, 13, 75.8681062153487),
  @(44481, 0, 9, 52.52407353370295),
  @(44482, 0, 9, 52.52407353370295),
  @(44483, 1, 6, 35.01604902246864),
  @(44484, 0, 3, 17.50802451123432),
  @(44485, 0, 3, 17.50802451123432),
  @(44486, 4, 7, 40.85205719288007),
  @(44487, 0, 5, 29.18004085205719),
  @(44488, 0, 5, 29.18004085205719),
  @(44489, 0, 5, 29.18004085205719),
  @(44490, 0, 4, 23.34403268164575),
  @(44491, 0, 4, 23.34403268164575),
  @(44492, 1, 5, 29.18004085205719),
  @(44493, 1, 2, 11.67201634082288),
  @(44494, 1, 3, 17.50802451123432),
  @(44495, 0, 3, 17.50802451123432),
  @(44496, 0, 3, 17.50802451123432),
  @(44497, 1, 4, 23.34403268164575),
  @(44498, 0, 4, 23.34403268164575),
  @(44499, 0, 3, 17.50802451123432),
  @(44500, 0, 2, 11.67201634082288),
  @(44501, 0, 1, 5.836008170411438),
  @(44502, 0, 1, 5.836008170411438),
  @(44503, 0, 1, 5.836008170411438),
  @(44504, 0, 0, 0),
  @(44505, 0, 0, 0),
  @(44506, 0, 0, 0),
  @(44507, 0, 0, 0),
  @(44508, 0, 0, 0),
  @(44509, 0, 0, 0),
  @(44510, 0, 0, 0),
  @(44511, 2, 2, 11.67201634082288),
  @(44512, 1, 3, 17.50802451123432),
  @(44513, 0, 3, 17.50802451123432),
  @(44514, 1, 4, 23.34403268164575),
  @(44515, 0, 4, 23.34403268164575),
  @(44516, 15, 19, 110.8841552378173),
  @(44517, 0, 19, 110.8841552378173),
  @(44518, 4, 21, 122.5561715786402),
  @(44519, 6, 26, 151.7362124306974),
  @(44520, 0, 26, 151.7362124306974),
  @(44521, 1, 26, 151.7362124306974),
  @(44522, 10, 36, 210.0962941348118),
  @(44523, 2, 23, 134.2281879194631),
  @(44524, 13, 36, 210.0962941348118),
  @(44525, 6, 38, 221.7683104756346),
  @(44526, 5, 37, 215.9323023052232),
  @(44527, 4, 41, 239.276334986869),
  @(44528, 4, 44, 256.7843594981033),
  @(44529, 5, 39, 227.6043186460461),
  @(44530, 2, 39, 227.6043186460461),
  @(44531, 0, 26, 151.7362124306974),
  @(44532, 9, 29, 169.2442369419317),
  @(44533, 1, 25, 145.900204260286),
  @(44534, 5, 26, 151.7362124306974),
  @(44535, 3, 25, 145.900204260286),
  @(44536, 2, 22, 128.3921797490516),
  @(44537, 2, 22, 128.3921797490516),
  @(44538, 3, 25, 145.900204260286)
)

$startRow = 386
$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
$endRow = $r - 1

# Column A carries the date number format + border/font style (same as the
# existing rows above) -- copy it down from the last pre-existing row so the
# new cells pick up style index 2 instead of getting a brand-new xf.
$ws.Range("A385").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0
